# Update the cached display text of every "datetimeFigureOut" date
# placeholder (Slide Master, every slide Layout, and the Notes Master)
# from "2/17/2018" to "4/5/2019" -- same as refreshing the date shown
# via Insert > Header & Footer > Apply to All, just pinned to a fixed
# date instead of "today".

$p = $ppt.ActivePresentation
$oldDate = "2/17/2018"
$newDate = "4/5/2019"
$ppPlaceholderDate = 16

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout hanging off the (single) slide master
$master = $p.SlideMaster
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateShape $layout.Shapes
}

# Notes Master (HasNotesMaster is unreliable in this host, so just
# reach for NotesMaster directly -- every presentation has one).
# The Notes Master's date placeholder shape doesn't take a direct
# TextRange.Text edit in this host, so go through the Header/Footer
# dialog's DateAndTime text instead (same field, different door).
# Note: DateAndTime.Text reads back empty in this host even though the
# write does land, so set it unconditionally rather than gating on the
# old value.
$nm = $p.NotesMaster
$nmDate = $nm.HeadersFooters.DateAndTime
$nmDate.Text = $newDate
